$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 198
$ws.Range("F6").Value = 847
$ws.Range("F7").Value = 4264
$ws.Range("F11").Value = 6267
$ws.Range("F13").Value = 471
$ws.Range("F14").Value = 2400
$ws.Range("F17").Value = 503
$ws.Range("F18").Value = 9
$ws.Range("F19").Value = 9409
$ws.Range("F21").Value = 2540
$ws.Range("F23").Value = 2344
$ws.Range("F24").Value = 2515
$ws.Range("F25").Value = 1398
$ws.Range("F26").Value = 253
$ws.Range("F27").Value = 1997
$ws.Range("F29").Value = 66
$ws.Range("F30").Value = 342
$ws.Range("F35").Value = 100
$ws.Range("F37").Value = 1241
$ws.Range("F38").Value = 1227
$ws.Range("F39").Value = 80
$ws.Range("F40").Value = 109
$ws.Range("F41").Value = 244
$ws.Range("F42").Value = 1586
$ws.Range("F43").Value = 2616
$ws.Range("F44").Value = 941
$ws.Range("F45").Value = 324
$ws.Range("F46").Value = 1261

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 921

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 198
$ws.Range("F4").Value = 921
$ws.Range("F9").Value = 4264
$ws.Range("F10").Value = 4264
$ws.Range("F14").Value = 6267
$ws.Range("F16").Value = 2400
$ws.Range("F18").Value = 503
$ws.Range("F19").Value = 9
$ws.Range("F20").Value = 9409
$ws.Range("F23").Value = 2540
$ws.Range("F24").Value = 2344
$ws.Range("F25").Value = 2515
$ws.Range("F26").Value = 253
$ws.Range("F27").Value = 1997
$ws.Range("F29").Value = 66
$ws.Range("F30").Value = 342
$ws.Range("F34").Value = 100
$ws.Range("F36").Value = 1241
$ws.Range("F37").Value = 1227
$ws.Range("F38").Value = 80
$ws.Range("F39").Value = 109
$ws.Range("F40").Value = 1586
$ws.Range("F41").Value = 2616
$ws.Range("F42").Value = 941
$ws.Range("F43").Value = 324
$ws.Range("F47").Value = 1261
